$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'247.25"
$ws.Range("G2").Value = "'10"

# Row 3
$ws.Range("D3").Value = "'21.86"
$ws.Range("G3").Value = "'10"

# Row 4
$ws.Range("D4").Value = "'5.400"
$ws.Range("G4").Value = "'10"

# Row 5
$ws.Range("D5").Value = "'0.05634"
$ws.Range("G5").Value = "'10"

# Row 6
$ws.Range("D6").Value = "'3.430"
$ws.Range("G6").Value = "'10"

# Row 7
$ws.Range("D7").Value = "'6.368"
$ws.Range("G7").Value = "'10"

# Row 8
$ws.Range("D8").Value = "'0.8170"
$ws.Range("G8").Value = "'10"

# Row 9
$ws.Range("D9").Value = "'0.9307"
$ws.Range("G9").Value = "'10"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1434"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "'10"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07518"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "'10"

# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03247"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'10"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03082"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "'10"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09323"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "'10"

# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.572"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "'10"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001593"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "'10"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04725"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "'10"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005783"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "'10"

# Row 19
$ws.Range("D19").Value = "'0.006391"
$ws.Range("G19").Value = "'10"

# Row 20
$ws.Range("D20").Value = "'0.005062"
$ws.Range("G20").Value = "'10"

# Row 21
$ws.Range("D21").Value = "'0.001034"
$ws.Range("G21").Value = "'10"

# Row 22
$ws.Range("G22").Value = "'10"

# Row 23
$ws.Range("D23").Value = "'3.747"
$ws.Range("G23").Value = "'10"

# Row 24
$ws.Range("D24").Value = "'2.179"
$ws.Range("G24").Value = "'10"

# Row 25
$ws.Range("D25").Value = "'0.3308"
$ws.Range("G25").Value = "'10"

# Row 26
$ws.Range("D26").Value = "'0.1320"
$ws.Range("G26").Value = "'10"

# Row 27
$ws.Range("E27").Value = "26AAXTokenAAB"
$ws.Range("G27").Value = "'10"

# Row 28
$ws.Range("D28").Value = "'0.0003001"
$ws.Range("G28").Value = "'10"

# Row 29
$ws.Range("G29").Value = "'10"

# Row 30
$ws.Range("G30").Value = "'10"

# Row 31
$ws.Range("G31").Value = "'10"

# Row 32
$ws.Range("G32").Value = "'10"

# Row 33
$ws.Range("G33").Value = "'10"

# Row 34
$ws.Range("G34").Value = "'10"

# Row 35
$ws.Range("G35").Value = "'10"

# Row 36
$ws.Range("G36").Value = "'10"

# Row 37
$ws.Range("G37").Value = "'10"

# Row 38
$ws.Range("G38").Value = "'10"

# Row 39
$ws.Range("G39").Value = "'10"

# Row 40
$ws.Range("G40").Value = "'10"

# Row 41
$ws.Range("D41").Value = "'0.006903"
$ws.Range("G41").Value = "'10"

# Row 42
$ws.Range("D42").Value = "'0.1064"
$ws.Range("G42").Value = "'10"

# Row 43
$ws.Range("G43").Value = "'10"

# Row 44
$ws.Range("D44").Value = "'0.008540"
$ws.Range("G44").Value = "'10"

# Row 45
$ws.Range("D45").Value = "'0.00005575"
$ws.Range("G45").Value = "'10"

# Row 46
$ws.Range("G46").Value = "'10"

# Row 47
$ws.Range("G47").Value = "'10"

# Row 48
$ws.Range("D48").Value = "'0.7804"
$ws.Range("G48").Value = "'10"

# Row 49
$ws.Range("D49").Value = "'0.1783"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("G49").Value = "'10"

# Row 50
$ws.Range("G50").Value = "'10"

# Row 51
$ws.Range("D51").Value = "'0.01010"
$ws.Range("G51").Value = "'10"
